$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# C10 was 18, restore it to 1 (per the commit's reverted/restored revision)
$ws.Range("C10").Value = 1
